# Weekly fruit/vegetable price update.
#
# The source table ("Hortaliza, Femacal de La Calera - Choclo") gains one
# new weekly record. It is inserted as the new row 703 (pushing the former
# rows 703-740 down to 704-741), so the sheet's used range grows from
# A1:R740 to A1:R741.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 703..740 down one row, opening up a blank row 703.
$ws.Rows.Item(703).Insert()

# Populate the newly-opened row 703 with the new weekly record.
$ws.Range("A703").Value = 3
$ws.Range("B703").Value = "Femacal de La Calera"
$ws.Range("C703").Value = "Coquimbo"
$ws.Range("D703").Value = 44753
$ws.Range("E703").Value = 5
$ws.Range("F703").Value = 100112024
$ws.Range("G703").Value = "Choclo"
$ws.Range("H703").Value = "Dulce o Americano"
$ws.Range("I703").Value = "Primera"
$ws.Range("J703").Value = 85
$ws.Range("K703").Value = 34000
$ws.Range("L703").Value = 35000
$ws.Range("M703").Value = 34471
$ws.Range("N703").Value = "$/malla 70 unidades"
$ws.Range("O703").Value = "Región de Arica y Parinacota"
$ws.Range("P703").Value = 492
$ws.Range("Q703").Value = 70
$ws.Range("R703").Value = "Hortaliza"
